# "Script to download new sets finished"
# Replace the placeholder single-set sheet ("Friday" + 7 sample cards) with
# the real, freshly-downloaded set header for Friday Night Magic 2001 (F01).
# The old sample card rows go away; the set has no cards logged yet, so
# row 2 is left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 becomes the full set title.
$ws.Range("A1").Value = "Friday Night Magic 2001 (F01)"

# Drop the old sample card rows (3-8); shift everything below them up.
$ws.Rows("3:8").Delete()

# No cards downloaded for this set yet - blank out row 2 (keeping it as an
# empty cell rather than removing it).
$ws.Range("A2").Value = " "
